$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update the F column sample numbers (F3:F27) from 1 to row-1 (2..26)
for ($row = 3; $row -le 27; $row++) {
    $ws.Cells.Item($row, 6).Value = $row - 1
}

# Update the view / selection to match the saved state (scroll position +
# the new active selection over the column that was just edited)
$win = $excel.ActiveWindow
$win.ScrollRow = 18
$win.ScrollColumn = 1
$ws.Range("F2:F27").Select()
